# Updated: st 29. 04. 2021
# Apply updated AgTests (F) and AgPosit (G) values for various rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{ Row = 271; F = 45668;  G = 1729 },
    @{ Row = 292; F = 82471 },
    @{ Row = 306; F = 75389 },
    @{ Row = 317; F = 63748 },
    @{ Row = 320; F = 71884;  G = 3313 },
    @{ Row = 332; F = 485563; G = 4810 },
    @{ Row = 335; F = 150424; G = 3764 },
    @{ Row = 338; F = 221553 },
    @{ Row = 348; F = 232776 },
    @{ Row = 349; F = 159419 },
    @{ Row = 350; F = 127017 },
    @{ Row = 353; F = 723586; G = 5291 },
    @{ Row = 357; F = 138232; G = 3013 },
    @{ Row = 360; F = 749755; G = 5140 },
    @{ Row = 363; G = 2750 },
    @{ Row = 393; F = 307452 },
    @{ Row = 395; F = 751195; G = 1960 },
    @{ Row = 399; F = 200886 },
    @{ Row = 401; F = 273431 },
    @{ Row = 408; F = 303699 },
    @{ Row = 409; F = 703249; G = 1004 },
    @{ Row = 410; F = 353881; G = 622 },
    @{ Row = 411; F = 224956; G = 824 },
    @{ Row = 412; F = 175739; G = 645 },
    @{ Row = 413; F = 148902 },
    @{ Row = 414; F = 145741; G = 550 },
    @{ Row = 415; F = 305435 },
    @{ Row = 416; F = 641063; G = 913 },
    @{ Row = 417; F = 325553; G = 556 },
    @{ Row = 418; F = 199523; G = 689 },
    @{ Row = 419; F = 145410; G = 499 }
)

foreach ($chg in $changes) {
    $r = $chg.Row
    if ($chg.ContainsKey('F')) {
        $ws.Cells.Item($r, 6).Value = $chg.F
    }
    if ($chg.ContainsKey('G')) {
        $ws.Cells.Item($r, 7).Value = $chg.G
    }
}
